$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$nm = $p.NotesMaster
Write-Output "SlideMaster.Theme:"
Write-Output $sm.Theme
Write-Output "NotesMaster.Theme:"
Write-Output $nm.Theme
try {
  $sm.Theme = $nm.Theme
  Write-Output "assignment worked"
} catch {
  Write-Output "assignment failed: $_"
}
